$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: new column H "LINKS" (same style as G1/"Notes") ---
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$ws.Range("H1").Value = "LINKS"

# Blue color used for every hyperlink font (VBA/COM color ints are BGR-ordered;
# 16711680 = 0x0000FF*256^2 -> renders as RGB FF0000FF).
$blue = 16711680

function Add-LinkCell {
    param(
        [string]$Cell,
        [string]$Text,
        [string]$Url
    )
    $ws.Range($Cell).Value = $Text
    $ws.Hyperlinks.Add($ws.Range($Cell), $Url) | Out-Null
    $ws.Range($Cell).Font.Underline = $true
    $ws.Range($Cell).Font.Color = $blue
}

# --- Row 2: OLED Module ---
Add-LinkCell "H2" "https://grobotronics.com/oled-module-0.91-128x32-i2c-white.html" "https://grobotronics.com/oled-module-0.91-128x32-i2c-white.html"

# --- Row 3: Rotary Encoder (display text has trailing newline) ---
Add-LinkCell "H3" "https://grobotronics.com/rotary-encoder-12mm-24p-r-with-switch.html`n" "https://grobotronics.com/rotary-encoder-12mm-24p-r-with-switch.html"
$ws.Rows.Item(3).AutoFit() | Out-Null

# --- Row 4: 74AHCT125 ---
Add-LinkCell "H4" "https://grobotronics.com/74ahct125-quad-level-shifter.html" "https://grobotronics.com/74ahct125-quad-level-shifter.html"

# --- Row 5: SK6812 (display text has trailing newline) ---
Add-LinkCell "H5" "https://grobotronics.com/led-8mm-rgb-programmable-ws2812-clone.html`n" "https://grobotronics.com/led-8mm-rgb-programmable-ws2812-clone.html"
$ws.Rows.Item(5).AutoFit() | Out-Null

# --- Row 6: 74HCT595D ---
Add-LinkCell "H6" "https://grobotronics.com/shift-register-8-bit-smd-74hct595d.html" "https://grobotronics.com/shift-register-8-bit-smd-74hct595d.html"

# --- Row 7: Capacitor 1000uF (display text has trailing newline) ---
Add-LinkCell "H7" "https://grobotronics.com/electrolytic-capacitor-16v-1000ufoem.html`n" "https://grobotronics.com/electrolytic-capacitor-16v-1000ufoem.html"
$ws.Rows.Item(7).AutoFit() | Out-Null

# --- Row 8: Πυκνωτής 100nf ---
Add-LinkCell "H8" "https://nettop.gr/index.php/hlektronika/capacitors/ceramic-capacitor/pyknotis-keramikos-104-100nf.html" "https://nettop.gr/index.php/hlektronika/capacitors/ceramic-capacitor/pyknotis-keramikos-104-100nf.html"

# --- Row 9: Raspberry Pi Pico ---
Add-LinkCell "H9" "https://nettop.gr/index.php/raspberry-pi/pico/raspberry-pi-pico.html?src=raspberrypi" "https://nettop.gr/index.php/raspberry-pi/pico/raspberry-pi-pico.html?src=raspberrypi"

# --- Row 10: 1N4148 ---
Add-LinkCell "H10" "https://nettop.gr/index.php/hlektronika/diode/diode-1n4148-small-signal-fast-switching-diode.html" "https://nettop.gr/index.php/hlektronika/diode/diode-1n4148-small-signal-fast-switching-diode.html"

# --- Row 11: Keychron K Pro Switches ---
Add-LinkCell "H11" "https://www.public.gr/product/gaming/pc-gaming/pc-gaming-accessories/diafora-pc-gaming-accessories/keychron-k-pro-switches-pliktrologiou-110-pack/2090270" "https://www.public.gr/product/gaming/pc-gaming/pc-gaming-accessories/diafora-pc-gaming-accessories/keychron-k-pro-switches-pliktrologiou-110-pack/2090270"

# --- Row 12: Glorious Stabilizers V2 ---
Add-LinkCell "H12" "https://shopflix.gr/p/SF-12868640/stabilizers-glorious-v2-gata1723-12843162" "https://shopflix.gr/p/SF-12868640/stabilizers-glorious-v2-gata1723-12843162"

# --- Row 13: Mountain Dolomite A Backlit Keycap Set ---
# Price correction: 5.22 -> 5.25 (shipping-adjusted running total updates too)
$ws.Range("D13").Value = 5.25
$ws.Range("E13").Value = 5.25
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F13").PasteSpecial(-4122) | Out-Null
$ws.Range("F13").Formula = "=5.25+F12+5"
Add-LinkCell "H13" "https://www.skroutz.gr/s/57409660/Mountain-Dolomite-A-Backlit-Keycap-Set.html" "https://www.skroutz.gr/s/57409660/Mountain-Dolomite-A-Backlit-Keycap-Set.html"

# --- Row 14: PCB (running total ripples from the row-13 correction) ---
$ws.Range("F14").Formula = "=F13+42.33"
Add-LinkCell "H14" "https://shopflix.gr/p/SF-200097568/mountain-dolomite-a-backlit-keycap-set" "https://shopflix.gr/p/SF-200097568/mountain-dolomite-a-backlit-keycap-set"

# --- Row 15: extra LINKS-only row for the JLCPCB site itself ---
Add-LinkCell "H15" "https://jlcpcb.com" "https://jlcpcb.com"
